$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 210912
$ws.Range("B8").Value = 120
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "ex 1.9-1.11 and material"

$ws.Range("A9").Select()
